# Revert "Remove deprecation, and instead only remove nanoDESI"
# i.e. re-add "NanoDESI" as an allowed assay_type value, before "NanoPOTS".

$wb = $excel.ActiveWorkbook

$listSheet = $wb.Worksheets.Item("assay_type list")

# Insert a new row above the existing row 1 (which currently holds "NanoPOTS"),
# so the existing value shifts down to row 2, and put "NanoDESI" into the new row 1.
$listSheet.Rows.Item(1).Insert() | Out-Null
$listSheet.Range("A1").Value = "NanoDESI"

# Update the data validation on the main sheet's assay_type column (L) so the
# list range covers both values and the error message reflects both options.
# Mutate the existing validation in place (rather than delete + re-add) so the
# dataValidations element order in the saved XML stays the same as before.
$mainSheet = $wb.Worksheets.Item("Export as TSV")
$col = $mainSheet.Range("L2:L1048576")

$col.Validation.Formula1 = "'assay_type list'!`$A`$1:`$A`$2"
$col.Validation.ErrorMessage = "Value must be one of: NanoDESI / NanoPOTS."
